$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new row at row 28 ("Ryze sim (470)") - shifts old rows 28-39
#    down to 29-40. Copy the formatting (styles + row height) from the row
#    that is about to land below it so the new row reuses existing style ids
#    instead of Excel minting brand-new (duplicate) style records.
# ---------------------------------------------------------------------------
$ws.Rows.Item(28).Insert()
$ws.Range("A29:H29").Copy()
$ws.Range("A28:H28").PasteSpecial(-4122)   # xlPasteFormats
$ws.Rows.Item(28).RowHeight = $ws.Rows.Item(29).RowHeight
$excel.CutCopyMode = 0

# Fill the new row's content
$ws.Range("B28").Value = "Ryze sim  (470)"
$ws.Range("C28").Value = 20
$ws.Range("G28").Formula = "=C28-D28+E28+F28"
$ws.Range("H28").Formula = "=G28*470"

# ---------------------------------------------------------------------------
# 2) Header banner text / date stamps
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Mangrove Communication  27.03.2025"
$ws.Range("A10").Value = "DAILY STOCK                         (27/03/2025) "

# ---------------------------------------------------------------------------
# 3) Sale & stock report table (rows 3-6)
# ---------------------------------------------------------------------------
$ws.Range("C3").Value = 243226
$ws.Range("D3").ClearContents()

$ws.Range("C4").Value = 215958
$ws.Range("D4").ClearContents()
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 7

$ws.Range("C5").Value = 156224
$ws.Range("D5").ClearContents()

$ws.Range("C6").Value = 240863
$ws.Range("D6").ClearContents()

# ---------------------------------------------------------------------------
# 4) I-top-up / balance table
# ---------------------------------------------------------------------------
$ws.Range("C13").Value = 28562

$ws.Range("C14").Value = 10322
$ws.Range("D14").Value = 856271
$ws.Range("E14").Value = 1166753

$ws.Range("C20").Value = 4000
$ws.Range("D20").ClearContents()

$ws.Range("C21").ClearContents()
$ws.Range("D21").ClearContents()

$ws.Range("C22").ClearContents()
$ws.Range("D22").ClearContents()

$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 1

$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 7

$ws.Range("C27").Value = 72

# ---------------------------------------------------------------------------
# 5) Cash / bank / expense summary (rows shifted down by the insert above;
#    use the NEW row numbers: old 33-39 -> new 34-40)
# ---------------------------------------------------------------------------
$ws.Range("H34").Value = 322679   # Cash
$ws.Range("H35").Value = 14670    # Bank
$ws.Range("H38").ClearContents()  # BG (was 150000, now blank)

# New "Loan to E-life" line, inserted right before GTAND TOTAL
$ws.Range("F39").Value = "Loan to E-life"
$ws.Range("G39").Style = $ws.Range("G40").Style
$ws.Range("H39").Value = 102000
$ws.Range("F39:G39").Merge()

# Grand-total formula now lists the rows in ascending order and includes H39
$ws.Range("H40").Formula = "=H33+H34+H35+H36+H37+H38+H39"

# ---------------------------------------------------------------------------
# 6) View state
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("K28").Select()
